$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the name in cell E7 from "Thu" to "Thuu"
$ws.Range("E7").Value = "Thuu"

# Reflect the final cell selection as in the edited file
$ws.Range("E7").Select()
